$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 100-101; existing rows 100-117 shift down to 102-119.
$ws.Rows("100:101").Insert()

# Shared/static values for this data set (same on every data row in this sheet).
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108003
$categoria   = "Maracuyá"
$variedad    = "Sin especificar"
$unidad      = "$/caja 20 kilos"
$origen      = "Región de Arica y Parinacota"
$kgUnidad    = 20

function Set-Fila($fila, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($fila, 1).Value = 1
    $ws.Cells.Item($fila, 2).Value = $mercado
    $ws.Cells.Item($fila, 3).Value = $region
    $ws.Cells.Item($fila, 4).Value = $fecha
    $ws.Cells.Item($fila, 5).Value = $codreg
    $ws.Cells.Item($fila, 6).Value = $tipo
    $ws.Cells.Item($fila, 7).Value = $productoId
    $ws.Cells.Item($fila, 8).Value = $producto
    $ws.Cells.Item($fila, 9).Value = $categoriaId
    $ws.Cells.Item($fila, 10).Value = $categoria
    $ws.Cells.Item($fila, 11).Value = $variedad
    $ws.Cells.Item($fila, 12).Value = $calidad
    $ws.Cells.Item($fila, 13).Value = $volumen
    $ws.Cells.Item($fila, 14).Value = $precioMin
    $ws.Cells.Item($fila, 15).Value = $precioMax
    $ws.Cells.Item($fila, 16).Value = $precioProm
    $ws.Cells.Item($fila, 17).Value = $unidad
    $ws.Cells.Item($fila, 18).Value = $origen
    $ws.Cells.Item($fila, 19).Value = $precioKg
    $ws.Cells.Item($fila, 20).Value = $kgUnidad
}

Set-Fila 100 44694 "Especial" 120 25000 26000 25500 1275
Set-Fila 101 44694 "Primera"  120 23000 24000 23500 1175
